$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "fra" language block (rows 6-9) with "spa"
$ws.Range("A6:A9").Value = "spa"

# Remove the remaining language blocks (ara, kan, hin, tam - rows 10-25)
$ws.Range("A10:D25").EntireRow.Delete()

# Update selection/view to match the post-edit state
$ws.Range("A6:A9").Select()
